# DescData.xlsx -- "unify the conception of DataNode, DataTable, Entity."
#
# The only workbook-semantic change in this commit is renaming the lone
# worksheet from "Property1" to "DataNode" (matching the new shared
# DataNode/DataTable/Entity naming scheme used across the other config
# sheets). The re-save also nudges the two header rows a touch shorter and
# leaves the cursor on D37 -- both are carried over here too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet: Property1 -> DataNode
$ws.Name = "DataNode"

# Header rows (1 and 8) shrink slightly: ht 28 -> 27
$ws.Rows.Item(1).RowHeight = 27
$ws.Rows.Item(8).RowHeight = 27

# Leave the live selection on D37, where editing ended
$ws.Activate() | Out-Null
$ws.Range("D37").Select() | Out-Null
